# Commit: Tue, Jun 16, 2020  6:05:43 PM
#
# The only semantically-meaningful, user-visible edit in the source diff is
# a change to the table style applied to the single table in the deck
# (the "Total Outflow / Total Inflow" plenary table on slide 16): its
# <a:tableStyleId> changes from the deck's custom style
# {2F563032-9014-4C15-AB7E-B1B76AF4B360} to the built-in PowerPoint table
# style {A65DC1F3-B025-406B-84C3-34F299455BB3}.
#
# (The rest of the diff is just ppt/theme/theme1.xml and ppt/theme/theme2.xml
# trading places/content in the package - the slide master's theme and the
# notes master's theme swap identities - which is not something that can be
# driven from the PowerPoint object model; it carries no visible change to
# any shape/text/formatting reachable via COM automation, so there is no
# Shape/Table/TextRange call that corresponds to it.)

$p = $ppt.ActivePresentation

$oldStyleId = "{2F563032-9014-4C15-AB7E-B1B76AF4B360}"
$newStyleId = "{A65DC1F3-B025-406B-84C3-34F299455BB3}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $tbl = $shape.Table
            if ($tbl.StyleId -eq $oldStyleId) {
                $tbl.ApplyStyle($newStyleId)
            }
        }
    }
}
